$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sprite_path values to point to the new Building subfolder
$ws.Range("E2").Value = "res://Asset/Building/shelf.png"
$ws.Range("E3").Value = "res://Asset/Building/stand.png"
$ws.Range("E4").Value = "res://Asset/Building/hang.png"

# Update column G width (col index 7) and drop bestFit/autofit sizing
$ws.Columns.Item(7).ColumnWidth = 10.83

# Update selection to E4
$ws.Range("E4").Select() | Out-Null
